$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 8 new rows at row 2 (pushes existing rows 2-21 down to rows 10-29)
$ws.Range("A2:A9").EntireRow.Insert()
# The insert copies formatting down from the header row; reset the new rows to the default (unstyled) look
$ws.Range("A2:H9").ClearFormats()

# Step 2: populate the 8 newly inserted rows (2-9) with new sensor data
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "falling"
$ws.Range("C2").Value = -0.5779368877410893
$ws.Range("D2").Value = 1.070879459381104
$ws.Range("E2").Value = 0.1698004633188247
$ws.Range("F2").Value = 0.0106901414692401
$ws.Range("G2").Value = -0.00335975876078
$ws.Range("H2").Value = 0.0360410511493682

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "falling"
$ws.Range("C3").Value = -0.6250030517578123
$ws.Range("D3").Value = 1.073733139038086
$ws.Range("E3").Value = 0.1257202506065367
$ws.Range("F3").Value = -0.007177666760981
$ws.Range("G3").Value = -0.0487165041267871
$ws.Range("H3").Value = 0.0716239511966705

$ws.Range("A4").Value = 200
$ws.Range("B4").Value = "falling"
$ws.Range("C4").Value = -0.5019012451171875
$ws.Range("D4").Value = 1.114973473548889
$ws.Range("E4").Value = 0.08085805475711817
$ws.Range("F4").Value = -0.0007635815418325
$ws.Range("G4").Value = -0.0448985956609249
$ws.Range("H4").Value = 0.0595593601465225

$ws.Range("A5").Value = 300
$ws.Range("B5").Value = "falling"
$ws.Range("C5").Value = -0.5343909263610841
$ws.Range("D5").Value = 1.139204859733582
$ws.Range("E5").Value = 0.1443376690149308
$ws.Range("F5").Value = 0.09666942805051799
$ws.Range("G5").Value = 0.0059559359215199
$ws.Range("H5").Value = 0.0488692186772823

$ws.Range("A6").Value = 400
$ws.Range("B6").Value = "falling"
$ws.Range("C6").Value = -0.5579452037811278
$ws.Range("D6").Value = 1.112600553035736
$ws.Range("E6").Value = 0.2124309107661247
$ws.Range("F6").Value = 0.1411098688840866
$ws.Range("G6").Value = 0.2434297949075698
$ws.Range("H6").Value = -0.0125227374956011

$ws.Range("A7").Value = 500
$ws.Range("B7").Value = "falling"
$ws.Range("C7").Value = -0.4796955108642578
$ws.Range("D7").Value = 1.016827774047851
$ws.Range("E7").Value = 0.1028751075267787
$ws.Range("F7").Value = -0.0161879286170005
$ws.Range("G7").Value = 0.07849618047475811
$ws.Range("H7").Value = 0.0746782794594764

$ws.Range("A8").Value = 600
$ws.Range("B8").Value = "falling"
$ws.Range("C8").Value = -0.5379581451416018
$ws.Range("D8").Value = 0.9855325698852542
$ws.Range("E8").Value = -0.2731702357530603
$ws.Range("F8").Value = 0.0678060427308082
$ws.Range("G8").Value = -0.026419922709465
$ws.Range("H8").Value = 0.0497855171561241

$ws.Range("A9").Value = 700
$ws.Range("B9").Value = "falling"
$ws.Range("C9").Value = -0.6476110458374021
$ws.Range("D9").Value = 1.080279231071473
$ws.Range("E9").Value = -0.8854551434516924
$ws.Range("F9").Value = -0.0740674138069152
$ws.Range("G9").Value = -0.4193589985370636
$ws.Range("H9").Value = 0.0155770638957619

# Step 3: append 2 new rows (30-31) with new sensor data after the existing data block
$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "falling"
$ws.Range("C30").Value = 0.00381779670715305
$ws.Range("D30").Value = 1.210070580244063
$ws.Range("E30").Value = -0.2140652965754264
$ws.Range("F30").Value = -0.0224492978304624
$ws.Range("G30").Value = 0.0058032199740409
$ws.Range("H30").Value = 0.0675006061792373

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "falling"
$ws.Range("C31").Value = -0.03790302276611289
$ws.Range("D31").Value = 1.083934617042542
$ws.Range("E31").Value = -0.04758519232273038
$ws.Range("F31").Value = -0.009010262787342
$ws.Range("G31").Value = -0.0429132841527462
$ws.Range("H31").Value = -0.06276640295982359

